$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value2 = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '30.548.48'
Set-TextValue $ws.Range('E2') '  -0.11%  '

Set-TextValue $ws.Range('D3') '1.918.35'
Set-TextValue $ws.Range('E3') '  -0.35%  '

Set-TextValue $ws.Range('D4') '1.001'
Set-TextValue $ws.Range('E4') '  +0.02%  '

Set-TextValue $ws.Range('D5') '245.54'
Set-TextValue $ws.Range('E5') '  -0.43%  '

Set-TextValue $ws.Range('D7') '0.4831'
Set-TextValue $ws.Range('E7') '  +1.76%  '

Set-TextValue $ws.Range('D8') '0.2895'
Set-TextValue $ws.Range('E8') '  -0.99%  '

Set-TextValue $ws.Range('E9') '  -1.33%  '

Set-TextValue $ws.Range('D10') '111.54'
Set-TextValue $ws.Range('E10') '  +5.16%  '

Set-TextValue $ws.Range('D11') '18.94'
Set-TextValue $ws.Range('E11') '  +2.84%  '

Set-TextValue $ws.Range('D12') '1.909.79'
Set-TextValue $ws.Range('E12') '  -0.73%  '

Set-TextValue $ws.Range('D13') '0.07561'
Set-TextValue $ws.Range('E13') '  -2.12%  '

Set-TextValue $ws.Range('D14') '5.281'
Set-TextValue $ws.Range('E14') '  -1.41%  '

Set-TextValue $ws.Range('D15') '0.6679'
Set-TextValue $ws.Range('E15') '  -0.69%  '

Set-TextValue $ws.Range('D16') '295.19'
Set-TextValue $ws.Range('E16') '  +2.22%  '

Set-TextValue $ws.Range('D17') '30.542.43'
Set-TextValue $ws.Range('E17') '  -0.23%  '

Set-TextValue $ws.Range('B18') 'Avalanche'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D18') '13.00'
Set-TextValue $ws.Range('E18') '  -0.29%  '

Set-TextValue $ws.Range('B19') 'Dai'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Range('D19') '1.001'
Set-TextValue $ws.Range('E19') '  +0.11%  '

Set-TextValue $ws.Range('D20') '0.000007603'
Set-TextValue $ws.Range('E20') '  -0.61%  '

Set-TextValue $ws.Range('D21') '5.563'
Set-TextValue $ws.Range('E21') '  +0.90%  '

Set-TextValue $ws.Range('D22') '2.168.53'
Set-TextValue $ws.Range('E22') '  +0.03%  '

Set-TextValue $ws.Range('D23') '1.001'
Set-TextValue $ws.Range('E23') '  -0.03%  '

Set-TextValue $ws.Range('D24') '6.440'
Set-TextValue $ws.Range('E24') '  +2.23%  '

Set-TextValue $ws.Range('D25') '9.456'
Set-TextValue $ws.Range('E25') '  +0.35%  '

Set-TextValue $ws.Range('D26') '164.86'

Set-TextValue $ws.Range('E27') '  -2.60%  '

Set-TextValue $ws.Range('D28') '2.094'
Set-TextValue $ws.Range('E28') '  -1.86%  '

Set-TextValue $ws.Range('D29') '0.1066'
Set-TextValue $ws.Range('E29') '  -2.02%  '

Set-TextValue $ws.Range('D30') '1.440'
Set-TextValue $ws.Range('E30') '  +5.76%  '

Set-TextValue $ws.Range('D31') '4.134'

Set-TextValue $ws.Range('D32') '4.064'
Set-TextValue $ws.Range('E32') '  -0.18%  '

Set-TextValue $ws.Range('D33') '0.05005'
Set-TextValue $ws.Range('E33') '  -1.46%  '

Set-TextValue $ws.Range('D34') '0.7407'
Set-TextValue $ws.Range('E34') '  -0.28%  '

Set-TextValue $ws.Range('D35') '1.137'
Set-TextValue $ws.Range('E35') '  -1.93%  '

Set-TextValue $ws.Range('D36') '0.9998'
Set-TextValue $ws.Range('E36') '  +0.05%  '

Set-TextValue $ws.Range('D37') '2.722'
Set-TextValue $ws.Range('E37') '  -1.01%  '

Set-TextValue $ws.Range('D38') '0.02014'
Set-TextValue $ws.Range('E38') '  -3.41%  '

Set-TextValue $ws.Range('D39') '2.684'
Set-TextValue $ws.Range('E39') '  -0.41%  '

Set-TextValue $ws.Range('D40') '110.69'
Set-TextValue $ws.Range('E40') '  -0.35%  '

Set-TextValue $ws.Range('D41') '2.012'
Set-TextValue $ws.Range('E41') '  -2.74%  '

Set-TextValue $ws.Range('D42') '0.4420'
Set-TextValue $ws.Range('E42') '  +0.12%  '

Set-TextValue $ws.Range('D43') '0.8664'
Set-TextValue $ws.Range('E43') '  -1.61%  '

Set-TextValue $ws.Range('B44') 'Aave'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D44') '70.57'
Set-TextValue $ws.Range('E44') '  +4.40%  '

Set-TextValue $ws.Range('B45') 'FraxShare'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws.Range('D45') '5.824'
Set-TextValue $ws.Range('E45') '  -2.09%  '

Set-TextValue $ws.Range('D47') '7.223'
Set-TextValue $ws.Range('E47') '  -0.86%  '

Set-TextValue $ws.Range('D48') '48.46'
Set-TextValue $ws.Range('E48') '  +2.64%  '

Set-TextValue $ws.Range('D49') '9.204'
Set-TextValue $ws.Range('E49') '  -2.12%  '

Set-TextValue $ws.Range('D50') '0.1232'
Set-TextValue $ws.Range('E50') '  -0.17%  '

Set-TextValue $ws.Range('D51') '0.2530'
Set-TextValue $ws.Range('E51') '  -0.58%  '
